$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Name of Algo - adjust imputed values for KNN result data
$ws.Range("E3").Value = 12.802
$ws.Range("B9").Value = 6.484999999999999
$ws.Range("B18").Value = 6.351
$ws.Range("B20").Value = 6.667999999999999
